# Auto-generated edit script applying the F-column (interest count) refresh
# and sheet4 ('全部类型') row content resync described in the diff.
$wb = $excel.ActiveWorkbook

# --- 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 574
$ws.Range("F4").Value = 469
$ws.Range("F5").Value = 299
$ws.Range("F6").Value = 202
$ws.Range("F8").Value = 1226
$ws.Range("F9").Value = 569
$ws.Range("F10").Value = 308
$ws.Range("F11").Value = 3
$ws.Range("F12").Value = 125
$ws.Range("F14").Value = 5746
$ws.Range("F15").Value = 1782
$ws.Range("F16").Value = 4176
$ws.Range("F17").Value = 433
$ws.Range("F18").Value = 238
$ws.Range("F20").Value = 4870
$ws.Range("F21").Value = 6255
$ws.Range("F23").Value = 1059
$ws.Range("F25").Value = 3783
$ws.Range("F26").Value = 498
$ws.Range("F28").Value = 195
$ws.Range("F29").Value = 132
$ws.Range("F30").Value = 992
$ws.Range("F32").Value = 474
$ws.Range("F33").Value = 561
$ws.Range("F34").Value = 1605
$ws.Range("F36").Value = 1725
$ws.Range("F37").Value = 198
$ws.Range("F39").Value = 1143
$ws.Range("F40").Value = 1337
$ws.Range("F43").Value = 3424
$ws.Range("F44").Value = 134
$ws.Range("F45").Value = 290
$ws.Range("F46").Value = 413
$ws.Range("F47").Value = 6
$ws.Range("F48").Value = 17
$ws.Range("F49").Value = 3890

# --- 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 1209

# --- 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 3936

# --- 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3936
$ws.Range("C3").Value = "杭州·次元盛典1.0"
$ws.Range("D3").Value = "康候圣街99号 顺丰创新中心"
$ws.Range("E3").Value = "2024.06.15 10:00-06.16 17:00"
$ws.Range("F3").Value = 2649
$ws.Range("G3").Value = 68
$ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=83672"
$ws.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202405/zF8i2b201716880393893.jpeg"
$ws.Range("C4").Value = "杭州·第三届动漫迷城嘉年华·毕业泳池"
$ws.Range("D4").Value = "风情大道2555号 第一世界大酒店(宋城杭州乐园旅游区湘湖店)"
$ws.Range("E4").Value = "2024.06.15 10:00-06.15 17:00"
$ws.Range("F4").Value = 574
$ws.Range("G4").Value = 70
$ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=84338"
$ws.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202405/8PQU0HpT1717144146761.jpeg"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2024-06-16"
$ws.Range("C5").Value = "杭州·动漫视界COS盛典"
$ws.Range("E5").Value = "2024.06.16 09:00-06.16 17:00"
$ws.Range("F5").Value = 469
$ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=85106"
$ws.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202405/Gv4ND6zs1715916154246.jpeg"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "2024-06-22"
$ws.Range("C6").Value = "杭州·巅峰对决·排球少年ONLY"
$ws.Range("D6").Value = "十四号大街431号 江滨篮球馆"
$ws.Range("E6").Value = "2024.06.22 10:00-06.22 17:00"
$ws.Range("F6").Value = 299
$ws.Range("G6").Value = 60
$ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=85095"
$ws.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202404/3WtpRjjo1714197500930.jpeg"
$ws.Range("C7").Value = "杭州·时光代理人「惊喜节拍」LIVE｜2024音乐巡演"
$ws.Range("D7").Value = "新北街85号三层G2-302 杭州大麦66 LIVEHOUSE"
$ws.Range("E7").Value = "2024.06.22 20:00-06.22 22:00"
$ws.Range("F7").Value = 1209
$ws.Range("G7").Value = 319
$ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=85043"
$ws.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202404/3nufasTp1714404961103.jpeg"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2024-06-23"
$ws.Range("C8").Value = "杭州·【早鸟5折】中西合奏·再现经典《青城山下·千年等一回》传世国风跨界音乐会"
$ws.Range("D8").Value = "曙光路31号 浙江音乐厅"
$ws.Range("E8").Value = "2024.06.23 15:00-06.23 21:00"
$ws.Range("F8").Value = 7
$ws.Range("G8").Value = 90
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=84597"
$ws.Range("I8").Value = "//i2.hdslb.com/bfs/openplatform/202404/jNu5hjYv1713514034369.jpeg"
$ws.Range("C9").Value = "杭州·第二届白日梦次元动漫嘉年华"
$ws.Range("D9").Value = "康候圣街99号 顺丰创新中心"
$ws.Range("E9").Value = "2024.06.23 10:00-06.23 17:00"
$ws.Range("F9").Value = 202
$ws.Range("G9").Value = 68
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=86307"
$ws.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202405/qHcyIUL31715752173541.jpeg"
$ws.Range("C10").Value = "杭州·第五人格ONLY2.0"
$ws.Range("D10").Value = "十四号大街431号 江滨篮球馆"
$ws.Range("F10").Value = 484
$ws.Range("G10").Value = 60
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=85710"
$ws.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202405/ULUN091G1715762966375.jpeg"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "2024-06-29"
$ws.Range("C11").Value = "杭州·《LALALAND爱乐之城》浪漫主题音乐会"
$ws.Range("D11").Value = "曙光路31号 浙江音乐厅"
$ws.Range("E11").Value = "2024.06.29 19:30-06.29 21:00"
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 100
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=87209"
$ws.Range("I11").Value = "//i1.hdslb.com/bfs/openplatform/202406/5s68NYdO1718171080184.jpeg"
$ws.Range("C12").Value = "杭州·乌托邦次元聚会3.0·二次元全女性夜场"
$ws.Range("D12").Value = "保淑路2号 The Queen皇后"
$ws.Range("E12").Value = "2024.06.29 13:00-06.29 19:00"
$ws.Range("F12").Value = 1226
$ws.Range("G12").Value = 188
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=84558"
$ws.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202404/XyOkWYv31713435061841.jpeg"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "2024-06-30"
$ws.Range("C13").Value = "杭州·热血番ONLY1.0"
$ws.Range("D13").Value = "康候圣街99号 顺丰创新中心"
$ws.Range("E13").Value = "2024.06.30 10:00-06.30 17:00"
$ws.Range("F13").Value = 569
$ws.Range("G13").Value = 68
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=85042"
$ws.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202406/ckVVTuNj1717752114555.jpeg"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "2024-07-04"
$ws.Range("C14").Value = "杭州·乙游Camerata杭州2.0"
$ws.Range("D14").Value = "杭海路601号江和美海洋广场1层 嘉宝丽酒店"
$ws.Range("E14").Value = "2024.07.04 10:00-07.04 17:00"
$ws.Range("F14").Value = 308
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=86777"
$ws.Range("I14").Value = "//i2.hdslb.com/bfs/openplatform/202406/aBDjuHlA1717403033570.jpeg"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "2024-07-06"
$ws.Range("C15").Value = "杭州·重逢·怀旧only"
$ws.Range("D15").Value = "丰庆路492号建冠龙禾商务中心A幢 杭州华礼宴国际礼宴中心(龙禾商务中心店)"
$ws.Range("E15").Value = "2024.07.06 09:00-07.06 17:00"
$ws.Range("F15").Value = 125
$ws.Range("G15").Value = 69
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=85742"
$ws.Range("I15").Value = "//i2.hdslb.com/bfs/openplatform/202405/qBeP0pEz1715399357252.jpeg"
$ws.Range("C16").Value = "杭州·黑执事only"
$ws.Range("D16").Value = "大岭山路156号 爱丽芬城堡"
$ws.Range("E16").Value = "2024.07.06 10:00-07.07 18:00"
$ws.Range("F16").Value = 359
$ws.Range("G16").Value = 160
$ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=86414"
$ws.Range("I16").Value = "//i1.hdslb.com/bfs/openplatform/202405/iP2cxk2w1716800288950.jpeg"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "2024-07-13"
$ws.Range("C17").Value = "杭州·【早鸟6折】《忱宴·渐渐被你吸引》热血动漫二次元ACG演唱会"
$ws.Range("D17").Value = "湖墅南路136-138号 浙话艺术剧院"
$ws.Range("E17").Value = "2024.07.13 19:30-07.13 21:30"
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 60
$ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=85011"
$ws.Range("I17").Value = "//i1.hdslb.com/bfs/openplatform/202404/2Gd8eLva1714379746993.jpeg"
$ws.Range("C18").Value = "杭州·代号鸢only-广陵大学"
$ws.Range("D18").Value = "康候圣街99号 顺丰创新中心"
$ws.Range("E18").Value = "2024.07.13 09:00-07.13 18:00"
$ws.Range("F18").Value = 1782
$ws.Range("G18").Value = 68
$ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=83289"
$ws.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202403/I3yffJ7Q1711344958258.png"
$ws.Range("C19").Value = "杭州·草莓动漫节"
$ws.Range("D19").Value = "中心路1号 白蓝地文创街区"
$ws.Range("E19").Value = "2024.07.13 09:00-07.14 17:00"
$ws.Range("F19").Value = 4176
$ws.Range("G19").Value = 70
$ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=84229"
$ws.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202405/yjYrwO301715760081303.jpeg"
$ws.Range("F20").Value = 4870
$ws.Range("F22").Value = 1059
$ws.Range("F24").Value = 3783
$ws.Range("F25").Value = 498
$ws.Range("F27").Value = 195
$ws.Range("F28").Value = 132
$ws.Range("F29").Value = 992
$ws.Range("F30").Value = 1414
$ws.Range("F31").Value = 474
$ws.Range("F32").Value = 561
$ws.Range("F34").Value = 1605
$ws.Range("F36").Value = 1725
$ws.Range("F43").Value = 3424
$ws.Range("F45").Value = 134
$ws.Range("F46").Value = 290
$ws.Range("F47").Value = 413
$ws.Range("F49").Value = 3890
